$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -5
$ws.Range("F3").Value = -13
$ws.Range("F5").Value = -1
$ws.Range("F6").Value = -3
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = -2
$ws.Range("F13").Value = -11
$ws.Range("F16").Value = -2
$ws.Range("F17").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("F21").Value = -1
$ws.Range("F32").Value = 3
$ws.Range("F41").Value = -3
